# This sheet is a weekly price log for "Berenjena" (eggplant) at the
# "Vega Modelo de Temuco" market. A new weekly record needs to be inserted
# at the top of the data block (row 194), pushing the existing rows
# 194-200 down to 195-201, and the new row 194 is populated with the
# latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 194; this shifts rows 194:200 down to
# 195:201 and keeps their formatting/styles intact (mirrors the dimension
# growing from A1:R200 to A1:R201).
$ws.Rows("194:194").Insert()

# Populate the newly inserted row 194 with the new weekly record.
$ws.Range("A194").Value = 10
$ws.Range("B194").Value = "Vega Modelo de Temuco"
$ws.Range("C194").Value = "La Araucanía"
$ws.Range("D194").Value = 44509
$ws.Range("E194").Value = 9
$ws.Range("F194").Value = 100112001
$ws.Range("G194").Value = "Berenjena"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 90
$ws.Range("K194").Value = 10000
$ws.Range("L194").Value = 11000
$ws.Range("M194").Value = 10556
$ws.Range("N194").Value = "$/caja 60 unidades"
$ws.Range("O194").Value = "Región de Arica y Parinacota"
$ws.Range("P194").Value = 176
$ws.Range("Q194").Value = 60
$ws.Range("R194").Value = "Hortaliza"

# Make sure the D194 date value keeps the date/time numeric style used by
# the rest of column D (in case Insert() did not propagate it).
$ws.Range("D194").NumberFormat = $ws.Range("D195").NumberFormat
